$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to remain
# text (matching the original inline-string storage) rather than being
# auto-converted to numeric values by Excel.

$ws.Range('D2').Value = '67.010.62'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').Value = '3.106.45'
$ws.Range('E3').Value = '  -0.11%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.26'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '177.93'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.59%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '3.103.29'
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('E9').Value = '  -1.61%  '
$ws.Range('E10').Value = '  -1.81%  '
$ws.Range('E11').Value = '  -1.75%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.467'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.25%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000241'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.21'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.61%  '
$ws.Range('E15').Value = '  -0.39%  '
$ws.Range('D16').Value = '3.624.06'
$ws.Range('E16').Value = '  -0.03%  '
$ws.Range('D17').Value = '66.960.42'
$ws.Range('E17').Value = '  -0.07%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.00'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.43%  '
$ws.Range('E19').Value = '  +2.15%  '
$ws.Range('D20').Value = '3.107.65'
$ws.Range('E20').Value = '  -0.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '485.15'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.31%  '
$ws.Range('E22').Value = '  -1.66%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.690'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.68'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.30%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.63'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.23'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.26'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.90%  '
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.05'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.80%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.28'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.59'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.59%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '28.04'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.18%  '
$ws.Range('E33').Value = '  -1.26%  '
$ws.Range('D34').Value = '0.0₃0940'
$ws.Range('E34').Value = '  -0.44%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '48.83'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.60'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.43%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.945'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.80%  '
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '49.17'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.66%  '
$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.309'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.05%  '
$ws.Range('E41').Value = '  -2.75%  '
$ws.Range('E42').Value = '  -0.16%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.30'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.89%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.68'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.30%  '
$ws.Range('D45').Value = '2.793.14'
$ws.Range('E45').Value = '  -0.76%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '372.21'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.47%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0344'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.71%  '
$ws.Range('E48').Value = '  -0.10%  '
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '24.94'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.43%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.23'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.01%  '
